$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 5274197
$ws.Cells.Item(17, 9).Value = 1299.5
$ws.Cells.Item(17, 10).Value = 5894537.5
$ws.Cells.Item(17, 11).Value = 3898.5
$ws.Cells.Item(17, 12).Value = 17683612.5
$ws.Cells.Item(17, 13).Value = -3730.5
$ws.Cells.Item(17, 14).Value = -17683948.5

$ws.Cells.Item(33, 8).Value = 163.78572
$ws.Cells.Item(33, 9).Value = 163.78572
$ws.Cells.Item(33, 11).Value = 163.78572
$ws.Cells.Item(33, 13).Value = 65.21428

$ws.Cells.Item(62, 8).Value = 4370.4136
$ws.Cells.Item(62, 9).Value = 3633.353
$ws.Cells.Item(62, 10).Value = 5414.5835
$ws.Cells.Item(62, 11).Value = 3633.353
$ws.Cells.Item(62, 12).Value = 5414.5835
$ws.Cells.Item(62, 13).Value = -3009.353
$ws.Cells.Item(62, 14).Value = -6662.5835

$ws.Cells.Item(65, 8).Value = 4370.4136
$ws.Cells.Item(65, 9).Value = 3633.353
$ws.Cells.Item(65, 10).Value = 5414.5835
$ws.Cells.Item(65, 11).Value = 18166.765
$ws.Cells.Item(65, 12).Value = 27072.9175
$ws.Cells.Item(65, 13).Value = -15046.765
$ws.Cells.Item(65, 14).Value = -33312.9175

$ws.Cells.Item(112, 8).Value = 3832478.2
$ws.Cells.Item(112, 10).Value = 3832478.2
$ws.Cells.Item(112, 12).Value = 11497434.6
$ws.Cells.Item(112, 14).Value = -11499650.6

$ws.Cells.Item(125, 8).Value = 720
$ws.Cells.Item(125, 9).Value = 600
$ws.Cells.Item(125, 10).Value = 800
$ws.Cells.Item(125, 11).Value = 5400
$ws.Cells.Item(125, 12).Value = 7200
$ws.Cells.Item(125, 13).Value = -2940
$ws.Cells.Item(125, 14).Value = -12120

$ws.Cells.Item(132, 8).Value = 1969.5686
$ws.Cells.Item(132, 9).Value = 2198.561
$ws.Cells.Item(132, 10).Value = 1030.7
$ws.Cells.Item(132, 11).Value = 6595.683000000001
$ws.Cells.Item(132, 12).Value = 3092.1
$ws.Cells.Item(132, 13).Value = -4065.683000000001
$ws.Cells.Item(132, 14).Value = -8152.1

$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()

$ws.Cells.Item(137, 8).Value = 1165.9318
$ws.Cells.Item(137, 9).Value = 1083.2667
$ws.Cells.Item(137, 11).Value = 3249.800099999999
$ws.Cells.Item(137, 13).Value = -699.8000999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 1000
$ws.Cells.Item(16, 9).Value = 1000
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 13).Value = -713

$ws.Cells.Item(32, 8).Value = 1862.1
$ws.Cells.Item(32, 9).Value = 1772.7579
$ws.Cells.Item(32, 10).Value = 3559.6
$ws.Cells.Item(32, 11).Value = 1772.7579
$ws.Cells.Item(32, 12).Value = 3559.6
$ws.Cells.Item(32, 13).Value = -1485.7579
$ws.Cells.Item(32, 14).Value = -4133.6

$ws.Cells.Item(61, 8).Value = 410246.9
$ws.Cells.Item(61, 9).Value = 419624.75
$ws.Cells.Item(61, 10).Value = 7000
$ws.Cells.Item(61, 11).Value = 419624.75
$ws.Cells.Item(61, 12).Value = 7000
$ws.Cells.Item(61, 13).Value = -419412.75
$ws.Cells.Item(61, 14).Value = -7424

$ws.Cells.Item(63, 8).Value = 2084406.2
$ws.Cells.Item(63, 9).Value = 1149.5
$ws.Cells.Item(63, 11).Value = 1149.5
$ws.Cells.Item(63, 13).Value = -463.5

$ws.Cells.Item(66, 8).Value = 2084406.2
$ws.Cells.Item(66, 9).Value = 1149.5
$ws.Cells.Item(66, 11).Value = 5747.5
$ws.Cells.Item(66, 13).Value = -2315.5

$ws.Cells.Item(74, 8).Value = 34485456
$ws.Cells.Item(74, 9).Value = 43480708
$ws.Cells.Item(74, 11).Value = 43480708
$ws.Cells.Item(74, 13).Value = -43479834

$ws.Cells.Item(77, 8).Value = 34485456
$ws.Cells.Item(77, 9).Value = 43480708
$ws.Cells.Item(77, 11).Value = 217403540
$ws.Cells.Item(77, 13).Value = -217399172

$ws.Cells.Item(122, 8).Value = 2105.6956
$ws.Cells.Item(122, 9).Value = 1467.7222
$ws.Cells.Item(122, 11).Value = 4403.1666
$ws.Cells.Item(122, 13).Value = -1953.1666

$ws.Cells.Item(136, 8).Value = 410246.9
$ws.Cells.Item(136, 9).Value = 419624.75
$ws.Cells.Item(136, 10).Value = 7000
$ws.Cells.Item(136, 11).Value = 1258874.25
$ws.Cells.Item(136, 12).Value = 21000
$ws.Cells.Item(136, 13).Value = -1256324.25
$ws.Cells.Item(136, 14).Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(25, 8).Value = 914
$ws.Cells.Item(25, 9).Value = 914
$ws.Cells.Item(25, 11).Value = 914
$ws.Cells.Item(25, 13).Value = -679

$ws.Cells.Item(134, 8).Value = 2427.7754
$ws.Cells.Item(134, 9).Value = 2477.413
$ws.Cells.Item(134, 10).Value = 1666.6666
$ws.Cells.Item(134, 11).Value = 7432.239
$ws.Cells.Item(134, 12).Value = 4999.9998
$ws.Cells.Item(134, 13).Value = -4897.239
$ws.Cells.Item(134, 14).Value = -10069.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3601.3044
$ws.Cells.Item(31, 9).Value = 2342.6667
$ws.Cells.Item(31, 11).Value = 2342.6667
$ws.Cells.Item(31, 13).Value = -2047.6667

$ws.Cells.Item(34, 8).Value = 3601.3044
$ws.Cells.Item(34, 9).Value = 2342.6667
$ws.Cells.Item(34, 11).Value = 2342.6667
$ws.Cells.Item(34, 13).Value = -2140.6667

$ws.Cells.Item(58, 8).Value = 11514.3545
$ws.Cells.Item(58, 9).Value = 904.6389
$ws.Cells.Item(58, 10).Value = 43343.5
$ws.Cells.Item(58, 11).Value = 904.6389
$ws.Cells.Item(58, 12).Value = 43343.5
$ws.Cells.Item(58, 13).Value = -701.6389
$ws.Cells.Item(58, 14).Value = -43749.5

$ws.Cells.Item(99, 8).Value = 18521538
$ws.Cells.Item(99, 9).Value = 2766.3125
$ws.Cells.Item(99, 10).Value = 45457936
$ws.Cells.Item(99, 11).Value = 2766.3125
$ws.Cells.Item(99, 12).Value = 45457936
$ws.Cells.Item(99, 13).Value = -1268.3125
$ws.Cells.Item(99, 14).Value = -45460932

$ws.Cells.Item(105, 8).Value = 6945785.5
$ws.Cells.Item(105, 9).Value = 7813802.5
$ws.Cells.Item(105, 11).Value = 7813802.5
$ws.Cells.Item(105, 13).Value = -7812055.5

$ws.Cells.Item(126, 8).Value = 18521538
$ws.Cells.Item(126, 9).Value = 2766.3125
$ws.Cells.Item(126, 10).Value = 45457936
$ws.Cells.Item(126, 11).Value = 8298.9375
$ws.Cells.Item(126, 12).Value = 136373808
$ws.Cells.Item(126, 13).Value = -5828.9375
$ws.Cells.Item(126, 14).Value = -136378748

$ws.Cells.Item(134, 8).Value = 643.4194
$ws.Cells.Item(134, 9).Value = 518.1111
$ws.Cells.Item(134, 10).Value = 1489.25
$ws.Cells.Item(134, 11).Value = 1554.3333
$ws.Cells.Item(134, 12).Value = 4467.75
$ws.Cells.Item(134, 13).Value = 980.6667000000002
$ws.Cells.Item(134, 14).Value = -9537.75

$ws.Cells.Item(136, 8).Value = 11514.3545
$ws.Cells.Item(136, 9).Value = 904.6389
$ws.Cells.Item(136, 10).Value = 43343.5
$ws.Cells.Item(136, 11).Value = 2713.9167
$ws.Cells.Item(136, 12).Value = 130030.5
$ws.Cells.Item(136, 13).Value = -163.9167000000002
$ws.Cells.Item(136, 14).Value = -135130.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(15, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 508.63635
$ws.Cells.Item(113, 10).Value = 549.375
$ws.Cells.Item(113, 12).Value = 1648.125
$ws.Cells.Item(113, 14).Value = -5988.125

$ws.Cells.Item(125, 8).Value = 5000
$ws.Cells.Item(125, 10).Value = 5000
$ws.Cells.Item(125, 12).Value = 15000
$ws.Cells.Item(125, 14).Value = -24840

$ws.Cells.Item(131, 8).Value = 784.24243
$ws.Cells.Item(131, 9).Value = 717.5
$ws.Cells.Item(131, 10).Value = 787.0526
$ws.Cells.Item(131, 11).Value = 2152.5
$ws.Cells.Item(131, 12).Value = 2361.1578
$ws.Cells.Item(131, 13).Value = 2887.5
$ws.Cells.Item(131, 14).Value = -12441.1578

$ws.Cells.Item(137, 8).Value = 15155090
$ws.Cells.Item(137, 9).Value = 466
$ws.Cells.Item(137, 10).Value = 27783944
$ws.Cells.Item(137, 11).Value = 1398
$ws.Cells.Item(137, 12).Value = 83351832
$ws.Cells.Item(137, 13).Value = 3702
$ws.Cells.Item(137, 14).Value = -83362032

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1745.8667
$ws.Cells.Item(100, 9).Value = 875.25
$ws.Cells.Item(100, 10).Value = 2740.8572
$ws.Cells.Item(100, 11).Value = 875.25
$ws.Cells.Item(100, 12).Value = 2740.8572
$ws.Cells.Item(100, 13).Value = -334.25
$ws.Cells.Item(100, 14).Value = -3822.8572

$ws.Cells.Item(132, 8).Value = 2573.4348
$ws.Cells.Item(132, 9).Value = 1557.25
$ws.Cells.Item(132, 10).Value = 4896.143
$ws.Cells.Item(132, 11).Value = 4671.75
$ws.Cells.Item(132, 12).Value = 14688.429
$ws.Cells.Item(132, 13).Value = -2141.75
$ws.Cells.Item(132, 14).Value = -19748.429

$ws.Cells.Item(136, 8).Value = 1286
$ws.Cells.Item(136, 9).Value = 1300.3
$ws.Cells.Item(136, 10).Value = 1000
$ws.Cells.Item(136, 11).Value = 3900.9
$ws.Cells.Item(136, 12).Value = 3000
$ws.Cells.Item(136, 13).Value = -1350.9
$ws.Cells.Item(136, 14).Value = -8100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 3497187.5
$ws.Cells.Item(107, 9).Value = 672.0909
$ws.Cells.Item(107, 10).Value = 22728022
$ws.Cells.Item(107, 11).Value = 2016.2727
$ws.Cells.Item(107, 12).Value = 68184066
$ws.Cells.Item(107, 13).Value = -96.27269999999999
$ws.Cells.Item(107, 14).Value = -68187906

$ws.Cells.Item(132, 8).Value = 791.5645
$ws.Cells.Item(132, 9).Value = 514.717
$ws.Cells.Item(132, 10).Value = 2421.889
$ws.Cells.Item(132, 11).Value = 1544.151
$ws.Cells.Item(132, 12).Value = 7265.667
$ws.Cells.Item(132, 13).Value = 985.8490000000002
$ws.Cells.Item(132, 14).Value = -12325.667
